# Round 1 clearing testcases + bug metrics update
# Adds a new manual test case (row 42) describing the "clearing of bid
# when number of bids exceed section size" scenario, and moves the
# active selection down to reflect the newly-added row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 42

$ws.Cells.Item($row, 1).Value = 41
$ws.Cells.Item($row, 2).Value = "Clearing of bid when number of bids exceed section size"
$ws.Cells.Item($row, 3).Value = "Bid for a course: round 1"
$ws.Cells.Item($row, 4).Value = "Admin logged in and presses Clear Round 1"
$ws.Cells.Item($row, 5).Value = "Load csv with 12 bids in IS100 section 1"
$ws.Cells.Item($row, 6).Value = "r1_test_section_size.zip file"
$ws.Cells.Item($row, 7).Value = "Sectionstudent table contains maggie.ng.2009 and neilson.ng.2009"
$ws.Cells.Item($row, 8).Value = "Sectionstudent contained all the bids who passed and failed"
$ws.Cells.Item($row, 9).Value = "Fail"

# Match the row height used by the other wrapped-text rows of similar length.
$ws.Rows.Item($row).RowHeight = 72.5

# Reflect the new active cell/selection that results from entering the row.
$ws.Range("F43").Select()
